$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Split the "How long is the average name? ... MAX7221" list paragraph
# (which carries the "_GoBack" bookmark right after its text) into two
# list paragraphs, inserting a new second bullet "How long are typical
# part numbers?" and moving the bookmark onto the end of that new bullet.
#
# We drive this off the bookmark's own Range so the edits naturally
# carry it along (inserting straight through Document.Range(pos,pos)
# at a "just before a paragraph mark" offset mis-anchors new bookmarks
# in this host, so we reuse/move the existing bookmark instead of
# deleting + re-adding it).
$goBack = $d.Bookmarks.Item("_GoBack")

# 1a. Insert a paragraph break right before the bookmark -- this turns
#     the bookmark (still collapsed) into the very start of a brand new
#     paragraph that inherits the same list formatting.
$breakRange = $goBack.Range.Duplicate
$breakRange.InsertBefore("`r")

# 1b. Insert the new bullet text right before the (now relocated)
#     bookmark, so the text lands in the new paragraph ahead of the
#     bookmark, which stays pinned to the end of that text.
$goBack = $d.Bookmarks.Item("_GoBack")
$textRange = $goBack.Range.Duplicate
$textRange.InsertBefore("How long are typical part numbers?")

# --- Edit 2 -----------------------------------------------------------
# Append a new paragraph after the "free text search" paragraph with the
# "Clicking on an applicable data point..." sentence.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# Re-fetch by index -- the Paragraph object returned by .Next() on a
# stale reference doesn't reliably re-seat after the structural edit,
# but re-querying the (now one-longer) collection does.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Clicking on an applicable data point (manufacturer, supplier, application) should bring up a selection window pre-loaded with a search on that metric. "
